$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: rename header and collapse round-status values (CLOSED/OPEN/COMPLETE)
#     into the new competition-level LOCKED/READY values.
#     (Column C is updated top-to-bottom first so the new shared strings
#      LOCKED / READY are introduced before UNLOCK / LOCK below.)
$ws.Range("C2").Value = "COMPETITION STATUS"
$ws.Range("C3").Value = "LOCKED"
$ws.Range("C4").Value = "LOCKED"
$ws.Range("C5").Value = "LOCKED"
$ws.Range("C6").Value = "READY"
$ws.Range("C7").Value = "READY"
$ws.Range("C8").Value = "LOCKED"
$ws.Range("C9").Value = "LOCKED"
$ws.Range("C10").Value = "LOCKED"
$ws.Range("C11").Value = "LOCKED"
$ws.Range("C12").Value = "LOCKED"

# --- Column B: rename the OPEN/CLOSE admin actions to UNLOCK/LOCK to match
#     the new competition-level locking actions.
$ws.Range("B6").Value = "UNLOCK"
$ws.Range("B9").Value = "LOCK"

# --- Remove the old "COMPETITION STATUS" block (rows 14-22: OPEN/CANCEL/DRAW/
#     WINNER + Admin can/Re-Open/Cancel) entirely - it's superseded by the
#     competition-level LOCKED/READY status now shown in column C above.
$ws.Range("B14:E22").EntireRow.Delete()

# --- Restore the selection to roughly where it was left after the edit.
$ws.Range("D18").Select() | Out-Null
